# Commit message: "Renaming files from kmos to kmcos"
#
# The document mentions the tool "kmos" in several places; rename every
# standalone occurrence of "kmos" to "kmcos" throughout the document body.
# (Note: "kmc_model" must NOT be touched - it does not contain "kmos".)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "kmos",    # FindText
    $true,     # MatchCase
    $false,    # MatchWholeWord
    $false,    # MatchWildcards
    $false,    # MatchSoundsLike
    $false,    # MatchAllWordForms
    $true,     # Forward
    1,         # Wrap (wdFindContinue)
    $false,    # Format
    "kmcos",   # ReplaceWith
    2          # Replace (wdReplaceAll)
)
